# Corrects the IFRS consolidated financial figures for 삼화페인트공업
# (rows 2-9, columns D:AJ) that were entered with wrong-magnitude values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5267
$ws.Range("E2").Value = 458
$ws.Range("F2").Value = 458
$ws.Range("G2").Value = 457
$ws.Range("H2").Value = 355
$ws.Range("I2").Value = 355
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5051
$ws.Range("L2").Value = 2254
$ws.Range("M2").Value = 2797
$ws.Range("N2").Value = 2795
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 112
$ws.Range("Q2").Value = 505
$ws.Range("R2").Value = -319
$ws.Range("S2").Value = -85
$ws.Range("T2").Value = 292
$ws.Range("U2").Value = 213
$ws.Range("V2").Value = 1067
$ws.Range("W2").Value = 8.699999999999999
$ws.Range("X2").Value = 6.75
$ws.Range("Y2").Value = 13.4
$ws.Range("Z2").Value = 7.25
$ws.Range("AA2").Value = 80.59999999999999
$ws.Range("AB2").Value = 2413.64
$ws.Range("AC2").Value = 1586
$ws.Range("AD2").Value = 9.359999999999999
$ws.Range("AE2").Value = 12897
$ws.Range("AF2").Value = 1.15
$ws.Range("AH2").Value = 2.69
$ws.Range("AI2").Value = 24.4
$ws.Range("AJ2").Value = 22400000

# Row 3
$ws.Range("D3").Value = 5072
$ws.Range("E3").Value = 317
$ws.Range("F3").Value = 317
$ws.Range("G3").Value = 295
$ws.Range("H3").Value = 248
$ws.Range("I3").Value = 248
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 5119
$ws.Range("L3").Value = 2241
$ws.Range("M3").Value = 2878
$ws.Range("N3").Value = 2875
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 113
$ws.Range("Q3").Value = 292
$ws.Range("R3").Value = -315
$ws.Range("S3").Value = -13
$ws.Range("T3").Value = 283
$ws.Range("U3").Value = 9
$ws.Range("V3").Value = 1153
$ws.Range("W3").Value = 6.25
$ws.Range("X3").Value = 4.88
$ws.Range("Y3").Value = 8.76
$ws.Range("Z3").Value = 4.87
$ws.Range("AA3").Value = 77.87
$ws.Range("AB3").Value = 2542.63
$ws.Range("AC3").Value = 1108
$ws.Range("AD3").Value = 11.46
$ws.Range("AE3").Value = 13556
$ws.Range("AF3").Value = 0.9399999999999999
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 3.94
$ws.Range("AI3").Value = 42.81
$ws.Range("AJ3").Value = 22668570

# Row 4
$ws.Range("D4").Value = 4822
$ws.Range("E4").Value = 189
$ws.Range("F4").Value = 189
$ws.Range("G4").Value = 191
$ws.Range("H4").Value = 137
$ws.Range("I4").Value = 137
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 5263
$ws.Range("L4").Value = 2306
$ws.Range("M4").Value = 2957
$ws.Range("N4").Value = 2950
$ws.Range("O4").Value = 7
$ws.Range("P4").Value = 127
$ws.Range("Q4").Value = 102
$ws.Range("R4").Value = -255
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = 233
$ws.Range("U4").Value = -131
$ws.Range("V4").Value = 1234
$ws.Range("W4").Value = 3.91
$ws.Range("X4").Value = 2.84
$ws.Range("Y4").Value = 4.72
$ws.Range("Z4").Value = 2.64
$ws.Range("AA4").Value = 77.97
$ws.Range("AB4").Value = 2387.65
$ws.Range("AC4").Value = 587
$ws.Range("AD4").Value = 16.88
$ws.Range("AE4").Value = 12695
$ws.Range("AF4").Value = 0.78
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 3.03
$ws.Range("AI4").Value = 50.74
$ws.Range("AJ4").Value = 25388674

# Row 5
$ws.Range("D5").Value = 4881
$ws.Range("E5").Value = 88
$ws.Range("F5").Value = 88
$ws.Range("G5").Value = 45
$ws.Range("H5").Value = 20
$ws.Range("I5").Value = 19
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5423
$ws.Range("L5").Value = 2499
$ws.Range("M5").Value = 2924
$ws.Range("N5").Value = 2918
$ws.Range("O5").Value = 7
$ws.Range("P5").Value = 131
$ws.Range("Q5").Value = 293
$ws.Range("R5").Value = -206
$ws.Range("S5").Value = -109
$ws.Range("T5").Value = 144
$ws.Range("U5").Value = 149
$ws.Range("V5").Value = 1148
$ws.Range("W5").Value = 1.8
$ws.Range("X5").Value = 0.4
$ws.Range("Y5").Value = 0.66
$ws.Range("Z5").Value = 0.37
$ws.Range("AA5").Value = 85.45999999999999
$ws.Range("AB5").Value = 2308.42
$ws.Range("AC5").Value = 75
$ws.Range("AD5").Value = 97.81999999999999
$ws.Range("AE5").Value = 12134
$ws.Range("AF5").Value = 0.61
$ws.Range("AG5").Value = 175
$ws.Range("AH5").Value = 2.37
$ws.Range("AI5").Value = 216.32
$ws.Range("AJ5").Value = 26196427

# Row 6
$ws.Range("D6").Value = 5242
$ws.Range("E6").Value = 79
$ws.Range("F6").Value = 79
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 8
$ws.Range("I6").Value = 8
$ws.Range("K6").Value = 5580
$ws.Range("L6").Value = 2762
$ws.Range("M6").Value = 2818
$ws.Range("N6").Value = 2809
$ws.Range("P6").Value = 132
$ws.Range("Q6").Value = 95
$ws.Range("R6").Value = -59
$ws.Range("S6").Value = 21
$ws.Range("T6").Value = 114
$ws.Range("U6").Value = -19
$ws.Range("V6").Value = 1354
$ws.Range("W6").Value = 1.5
$ws.Range("X6").Value = 0.15
$ws.Range("Y6").Value = 0.27
$ws.Range("Z6").Value = 0.14
$ws.Range("AA6").Value = 98.02
$ws.Range("AB6").Value = 2280.88
$ws.Range("AC6").Value = 29
$ws.Range("AD6").Value = 240.98
$ws.Range("AE6").Value = 12294
$ws.Range("AF6").Value = 0.5600000000000001
$ws.Range("AG6").Value = 125
$ws.Range("AH6").Value = 1.81
$ws.Range("AI6").Value = 377.13
$ws.Range("AJ6").Value = 26438751

# Row 7
$ws.Range("D7").Value = 5347
$ws.Range("E7").Value = 112
$ws.Range("G7").Value = 80
$ws.Range("H7").Value = 55
$ws.Range("I7").Value = 55
$ws.Range("K7").Value = 5662
$ws.Range("L7").Value = 2816
$ws.Range("M7").Value = 2846
$ws.Range("N7").Value = 2838
$ws.Range("P7").Value = 132
$ws.Range("Q7").Value = 267
$ws.Range("R7").Value = -56
$ws.Range("S7").Value = 4
$ws.Range("T7").Value = 68
$ws.Range("W7").Value = 2.1
$ws.Range("X7").Value = 1.03
$ws.Range("Y7").Value = 1.95
$ws.Range("Z7").Value = 0.98
$ws.Range("AA7").Value = 98.95
$ws.Range("AC7").Value = 208
$ws.Range("AD7").Value = 23.99
$ws.Range("AE7").Value = 12536
$ws.Range("AF7").Value = 0.4
$ws.Range("U7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 5507
$ws.Range("E8").Value = 138
$ws.Range("G8").Value = 97
$ws.Range("H8").Value = 74
$ws.Range("I8").Value = 74
$ws.Range("K8").Value = 5763
$ws.Range("L8").Value = 2870
$ws.Range("M8").Value = 2893
$ws.Range("N8").Value = 2885
$ws.Range("P8").Value = 132
$ws.Range("Q8").Value = 182
$ws.Range("R8").Value = -22
$ws.Range("S8").Value = -17
$ws.Range("T8").Value = 34
$ws.Range("W8").Value = 2.51
$ws.Range("X8").Value = 1.34
$ws.Range("Y8").Value = 2.59
$ws.Range("Z8").Value = 1.29
$ws.Range("AA8").Value = 99.2
$ws.Range("AC8").Value = 280
$ws.Range("AD8").Value = 17.35
$ws.Range("AE8").Value = 12744
$ws.Range("AF8").Value = 0.38
$ws.Range("U8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 5672
$ws.Range("E9").Value = 163
$ws.Range("G9").Value = 120
$ws.Range("H9").Value = 91
$ws.Range("I9").Value = 91
$ws.Range("K9").Value = 5806
$ws.Range("L9").Value = 2849
$ws.Range("M9").Value = 2958
$ws.Range("N9").Value = 2949
$ws.Range("P9").Value = 132
$ws.Range("Q9").Value = 190
$ws.Range("R9").Value = -3
$ws.Range("S9").Value = -94
$ws.Range("T9").Value = 15
$ws.Range("W9").Value = 2.87
$ws.Range("X9").Value = 1.6
$ws.Range("Y9").Value = 3.12
$ws.Range("Z9").Value = 1.57
$ws.Range("AA9").Value = 96.31999999999999
$ws.Range("AC9").Value = 344
$ws.Range("AD9").Value = 14.11
$ws.Range("AE9").Value = 13026
$ws.Range("AF9").Value = 0.37
$ws.Range("U9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

Write-Host "applied edits"
